$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276206556/TJaAJ2AJpb-camera-video-segment-driver-1743276209056.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSGRTBMXLY%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T070004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBwaCXVzLXdlc3QtMiJHMEUCIQC6qS59rLkPC17OmZgjRc%2BYyxhfX3YlVHis%2BJeZyxl02QIgJR8zAVTjyUvj0CmWm2bB7MfMnTfoWfDKr8YQZQciqqMq5gMIhf%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDGMZyK2QNyTMboqq3Sq6AyTWBDmEjvHSS7qDM9N9xdXIw%2BXirmg3UBdWSBH85PaULJiEN1%2B1gQY5uS4WWqVKEApioJVrDRCOnwG68ltDNvYxaTPfh27Ygz%2FSmITf6Soy%2FZ9%2Far311JEBAJ%2BD%2BhO%2BELXBH%2BhddsOj8DXo9YSX%2BdCXHzYePar%2F31yYToHaAobkr1xhX7DV644lIz8ocjPly%2BwFL8M2%2FGoyk4P6ZzYtW19b%2BmnSLLdsmFqYyLiYrE8zd7DQBtfsfRPY9oxlB%2BRzb8fUED3IOs0w9vOgMrSC%2Bx5UzFqHy8SLdXAv%2F1hQD1pOJbE%2FALOxtNvEXOuxGTOnWlIpomUem8W3KvjrIS%2BOyWsW%2BnyaUN%2FFf%2FG%2BTfYiuUizzsgFke58KTaIitG%2F8jklbjc0G4GoOJ2PP3mQKKLuSkV0BzDruDj%2BmN3%2Biw0hs7CkCfQ4aSmy2BFkEyWXL638Yy2fynTy39rzsZ4qob2dxulYTjp%2BOfph4FDEu8zTipdLxkMJWjtqOjXfRESfDyhGBRxUb9kXAIYRPJfiKdaU4EW%2BVkieKgVWdhFyVtQFJ9t3sIG46fRilYQIWkMIEgKCPMzYgInWOWqN7X0wu4mjvwY6pQHbJ0oE5%2Fkk4hVnQhrXWcuPqaVkGhuZG188waGCWK%2FMVS7FR3eMBtm62bIga8A7jx%2B0nnd2xOvZsW%2FPQXc9Z663uknPEG6YNBB45Bq6IWE2Ivn1%2FDedy9vvtyHa6bIMUXVjfIemVvmXIQYq27vUHZZaX2v9p1a8zGoCJUQ7%2BReVC2oU31r%2F91je4Zs%2F4ug4736F6Sp1o%2FTAa5XffMQiG6dRdslyCbs%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2015%3A00%3A04%20GMT&X-Amz-Signature=4e05c02e6d353acd85a745ac1a29105115cdec07e577b80cb6bf52d1c9626828'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276164863/OF6jgx1Lcp-camera-video-segment-driver-1743276167363.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSGRTBMXLY%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T070004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBwaCXVzLXdlc3QtMiJHMEUCIQC6qS59rLkPC17OmZgjRc%2BYyxhfX3YlVHis%2BJeZyxl02QIgJR8zAVTjyUvj0CmWm2bB7MfMnTfoWfDKr8YQZQciqqMq5gMIhf%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDGMZyK2QNyTMboqq3Sq6AyTWBDmEjvHSS7qDM9N9xdXIw%2BXirmg3UBdWSBH85PaULJiEN1%2B1gQY5uS4WWqVKEApioJVrDRCOnwG68ltDNvYxaTPfh27Ygz%2FSmITf6Soy%2FZ9%2Far311JEBAJ%2BD%2BhO%2BELXBH%2BhddsOj8DXo9YSX%2BdCXHzYePar%2F31yYToHaAobkr1xhX7DV644lIz8ocjPly%2BwFL8M2%2FGoyk4P6ZzYtW19b%2BmnSLLdsmFqYyLiYrE8zd7DQBtfsfRPY9oxlB%2BRzb8fUED3IOs0w9vOgMrSC%2Bx5UzFqHy8SLdXAv%2F1hQD1pOJbE%2FALOxtNvEXOuxGTOnWlIpomUem8W3KvjrIS%2BOyWsW%2BnyaUN%2FFf%2FG%2BTfYiuUizzsgFke58KTaIitG%2F8jklbjc0G4GoOJ2PP3mQKKLuSkV0BzDruDj%2BmN3%2Biw0hs7CkCfQ4aSmy2BFkEyWXL638Yy2fynTy39rzsZ4qob2dxulYTjp%2BOfph4FDEu8zTipdLxkMJWjtqOjXfRESfDyhGBRxUb9kXAIYRPJfiKdaU4EW%2BVkieKgVWdhFyVtQFJ9t3sIG46fRilYQIWkMIEgKCPMzYgInWOWqN7X0wu4mjvwY6pQHbJ0oE5%2Fkk4hVnQhrXWcuPqaVkGhuZG188waGCWK%2FMVS7FR3eMBtm62bIga8A7jx%2B0nnd2xOvZsW%2FPQXc9Z663uknPEG6YNBB45Bq6IWE2Ivn1%2FDedy9vvtyHa6bIMUXVjfIemVvmXIQYq27vUHZZaX2v9p1a8zGoCJUQ7%2BReVC2oU31r%2F91je4Zs%2F4ug4736F6Sp1o%2FTAa5XffMQiG6dRdslyCbs%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2015%3A00%3A04%20GMT&X-Amz-Signature=2ce54363a47f172fdc64507a51d8f92f80562797c065f5405ed44e7ee74541a0'
$ws.Range("K4").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276108738/gAxSKLP636-camera-video-segment-driver-1743276111238.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSGRTBMXLY%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T070004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBwaCXVzLXdlc3QtMiJHMEUCIQC6qS59rLkPC17OmZgjRc%2BYyxhfX3YlVHis%2BJeZyxl02QIgJR8zAVTjyUvj0CmWm2bB7MfMnTfoWfDKr8YQZQciqqMq5gMIhf%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDGMZyK2QNyTMboqq3Sq6AyTWBDmEjvHSS7qDM9N9xdXIw%2BXirmg3UBdWSBH85PaULJiEN1%2B1gQY5uS4WWqVKEApioJVrDRCOnwG68ltDNvYxaTPfh27Ygz%2FSmITf6Soy%2FZ9%2Far311JEBAJ%2BD%2BhO%2BELXBH%2BhddsOj8DXo9YSX%2BdCXHzYePar%2F31yYToHaAobkr1xhX7DV644lIz8ocjPly%2BwFL8M2%2FGoyk4P6ZzYtW19b%2BmnSLLdsmFqYyLiYrE8zd7DQBtfsfRPY9oxlB%2BRzb8fUED3IOs0w9vOgMrSC%2Bx5UzFqHy8SLdXAv%2F1hQD1pOJbE%2FALOxtNvEXOuxGTOnWlIpomUem8W3KvjrIS%2BOyWsW%2BnyaUN%2FFf%2FG%2BTfYiuUizzsgFke58KTaIitG%2F8jklbjc0G4GoOJ2PP3mQKKLuSkV0BzDruDj%2BmN3%2Biw0hs7CkCfQ4aSmy2BFkEyWXL638Yy2fynTy39rzsZ4qob2dxulYTjp%2BOfph4FDEu8zTipdLxkMJWjtqOjXfRESfDyhGBRxUb9kXAIYRPJfiKdaU4EW%2BVkieKgVWdhFyVtQFJ9t3sIG46fRilYQIWkMIEgKCPMzYgInWOWqN7X0wu4mjvwY6pQHbJ0oE5%2Fkk4hVnQhrXWcuPqaVkGhuZG188waGCWK%2FMVS7FR3eMBtm62bIga8A7jx%2B0nnd2xOvZsW%2FPQXc9Z663uknPEG6YNBB45Bq6IWE2Ivn1%2FDedy9vvtyHa6bIMUXVjfIemVvmXIQYq27vUHZZaX2v9p1a8zGoCJUQ7%2BReVC2oU31r%2F91je4Zs%2F4ug4736F6Sp1o%2FTAa5XffMQiG6dRdslyCbs%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2015%3A00%3A04%20GMT&X-Amz-Signature=7142c86fd3d83a3ed540b2ebd0584c028df848853430d5ed15a84fad2f975315'
$ws.Range("K5").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205341/1743265652740/hq4516Mzws-camera-video-segment-driver-1743265655240.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSGRTBMXLY%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T070003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBwaCXVzLXdlc3QtMiJHMEUCIQC6qS59rLkPC17OmZgjRc%2BYyxhfX3YlVHis%2BJeZyxl02QIgJR8zAVTjyUvj0CmWm2bB7MfMnTfoWfDKr8YQZQciqqMq5gMIhf%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDGMZyK2QNyTMboqq3Sq6AyTWBDmEjvHSS7qDM9N9xdXIw%2BXirmg3UBdWSBH85PaULJiEN1%2B1gQY5uS4WWqVKEApioJVrDRCOnwG68ltDNvYxaTPfh27Ygz%2FSmITf6Soy%2FZ9%2Far311JEBAJ%2BD%2BhO%2BELXBH%2BhddsOj8DXo9YSX%2BdCXHzYePar%2F31yYToHaAobkr1xhX7DV644lIz8ocjPly%2BwFL8M2%2FGoyk4P6ZzYtW19b%2BmnSLLdsmFqYyLiYrE8zd7DQBtfsfRPY9oxlB%2BRzb8fUED3IOs0w9vOgMrSC%2Bx5UzFqHy8SLdXAv%2F1hQD1pOJbE%2FALOxtNvEXOuxGTOnWlIpomUem8W3KvjrIS%2BOyWsW%2BnyaUN%2FFf%2FG%2BTfYiuUizzsgFke58KTaIitG%2F8jklbjc0G4GoOJ2PP3mQKKLuSkV0BzDruDj%2BmN3%2Biw0hs7CkCfQ4aSmy2BFkEyWXL638Yy2fynTy39rzsZ4qob2dxulYTjp%2BOfph4FDEu8zTipdLxkMJWjtqOjXfRESfDyhGBRxUb9kXAIYRPJfiKdaU4EW%2BVkieKgVWdhFyVtQFJ9t3sIG46fRilYQIWkMIEgKCPMzYgInWOWqN7X0wu4mjvwY6pQHbJ0oE5%2Fkk4hVnQhrXWcuPqaVkGhuZG188waGCWK%2FMVS7FR3eMBtm62bIga8A7jx%2B0nnd2xOvZsW%2FPQXc9Z663uknPEG6YNBB45Bq6IWE2Ivn1%2FDedy9vvtyHa6bIMUXVjfIemVvmXIQYq27vUHZZaX2v9p1a8zGoCJUQ7%2BReVC2oU31r%2F91je4Zs%2F4ug4736F6Sp1o%2FTAa5XffMQiG6dRdslyCbs%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2015%3A00%3A03%20GMT&X-Amz-Signature=4cfcaa97b5a2a79712899a52584c2863ad9b540da2c6834d2cbb507c1004cd5a'
$ws.Range("K6").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743255999969/Yjq7ZuVo4p-camera-video-segment-driver-1743256002469.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSGRTBMXLY%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T070003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBwaCXVzLXdlc3QtMiJHMEUCIQC6qS59rLkPC17OmZgjRc%2BYyxhfX3YlVHis%2BJeZyxl02QIgJR8zAVTjyUvj0CmWm2bB7MfMnTfoWfDKr8YQZQciqqMq5gMIhf%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDGMZyK2QNyTMboqq3Sq6AyTWBDmEjvHSS7qDM9N9xdXIw%2BXirmg3UBdWSBH85PaULJiEN1%2B1gQY5uS4WWqVKEApioJVrDRCOnwG68ltDNvYxaTPfh27Ygz%2FSmITf6Soy%2FZ9%2Far311JEBAJ%2BD%2BhO%2BELXBH%2BhddsOj8DXo9YSX%2BdCXHzYePar%2F31yYToHaAobkr1xhX7DV644lIz8ocjPly%2BwFL8M2%2FGoyk4P6ZzYtW19b%2BmnSLLdsmFqYyLiYrE8zd7DQBtfsfRPY9oxlB%2BRzb8fUED3IOs0w9vOgMrSC%2Bx5UzFqHy8SLdXAv%2F1hQD1pOJbE%2FALOxtNvEXOuxGTOnWlIpomUem8W3KvjrIS%2BOyWsW%2BnyaUN%2FFf%2FG%2BTfYiuUizzsgFke58KTaIitG%2F8jklbjc0G4GoOJ2PP3mQKKLuSkV0BzDruDj%2BmN3%2Biw0hs7CkCfQ4aSmy2BFkEyWXL638Yy2fynTy39rzsZ4qob2dxulYTjp%2BOfph4FDEu8zTipdLxkMJWjtqOjXfRESfDyhGBRxUb9kXAIYRPJfiKdaU4EW%2BVkieKgVWdhFyVtQFJ9t3sIG46fRilYQIWkMIEgKCPMzYgInWOWqN7X0wu4mjvwY6pQHbJ0oE5%2Fkk4hVnQhrXWcuPqaVkGhuZG188waGCWK%2FMVS7FR3eMBtm62bIga8A7jx%2B0nnd2xOvZsW%2FPQXc9Z663uknPEG6YNBB45Bq6IWE2Ivn1%2FDedy9vvtyHa6bIMUXVjfIemVvmXIQYq27vUHZZaX2v9p1a8zGoCJUQ7%2BReVC2oU31r%2F91je4Zs%2F4ug4736F6Sp1o%2FTAa5XffMQiG6dRdslyCbs%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2015%3A00%3A03%20GMT&X-Amz-Signature=71300b3e3538e5e1066486bf2455b6b55bcfe7362484bf7b77f5829e77506200'
